# Project 2 finalized w/ new due date!
# Update the "datetimeFigureOut" date placeholder (Insert > Header & Footer
# fixed date) from 4/2/2020 to 4/30/2020 across the slide master and every
# slide layout.

$p = $ppt.ActivePresentation

$ppPlaceholderDate = 16
$oldDate = "4/2/2020"
$newDate = "4/30/2020"

function Update-DatePlaceholder($shapes, $oldText, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $sh.TextFrame.TextRange.Text -eq $oldText) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes $oldDate $newDate

# Every slide layout off the master
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    Update-DatePlaceholder $layout.Shapes $oldDate $newDate
}
